$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add four new quota sheets (gender, age quota AU, CA, TR, UA) by
#    copying the existing "quotas_SK" sheet (same layout/formulas) and
#    then overwriting the country-specific figures.
# ---------------------------------------------------------------------

function Add-QuotaSheet {
    param(
        [string]$Name,
        [double]$FemaleShare,
        [double]$Age1,
        [double]$Age2,
        [double]$Age3,
        [double]$Age4,
        [double]$Age5
    )

    $template = $wb.Worksheets.Item("quotas_SK")
    $template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $ws = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws.Name = $Name

    # gender split: C2 holds the literal "man" share, B2 derives the
    # "woman" share as 1-C2 (rather than the ratio-based formula used by
    # the older template sheets)
    $ws.Range("C2").Value = $FemaleShare
    $ws.Range("B2").Formula = "=1-C2"
    $ws.Range("F3").ClearContents()

    # age bands
    $ws.Range("B8").Value = $Age1
    $ws.Range("C8").Value = $Age2
    $ws.Range("D8").Value = $Age3
    $ws.Range("E8").Value = $Age4
    $ws.Range("F8").Value = $Age5

    return $ws
}

Add-QuotaSheet -Name "quota_AU" `
    -FemaleShare 0.49354379369142626 `
    -Age1 0.11192255877789467 -Age2 0.18591237270266067 -Age3 0.26180562458671341 -Age4 0.23014443705193119 -Age5 0.21021500688079983 | Out-Null

Add-QuotaSheet -Name "quota_CA" `
    -FemaleShare 0.4927381777223736 `
    -Age1 0.10402815988780871 -Age2 0.17502163559091188 -Age3 0.24488811098766669 -Age4 0.25292269475436224 -Age5 0.22313939877924996 | Out-Null

Add-QuotaSheet -Name "quota_TR" `
    -FemaleShare 0.48657573802133475 `
    -Age1 0.1582325601298683 -Age2 0.2131137356790879 -Age3 0.29681998391047659 -Age4 0.20571250504618752 -Age5 0.12612121523437902 | Out-Null

$wsUA = Add-QuotaSheet -Name "quota_UA" `
    -FemaleShare 0.45142595728437557 `
    -Age1 8.2159805419933327E-2 -Age2 0.17834324173209759 -Age3 0.28227842928850105 -Age4 0.24861654266018471 -Age5 0.20860198089928325

# ---------------------------------------------------------------------
# 2. Restore a plain selection on each new sheet (Excel leaves the
#    inherited B8:F8 selection from the copied template otherwise).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("quota_AU").Range("C2").Select() | Out-Null
$wb.Worksheets.Item("quota_CA").Range("C2").Select() | Out-Null
$wb.Worksheets.Item("quota_TR").Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. View-state tweaks on pre-existing sheets.
# ---------------------------------------------------------------------

# Specificities: unfreeze/re-freeze at A2 (was A82) and move the
# selection to V4 (was R94)
$wsSpec = $wb.Worksheets.Item("Specificities")
$wsSpec.Activate()
$win = $excel.Windows.Item(1)
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.FreezePanes = $true
$wsSpec.Range("V4").Select() | Out-Null

# quotas_CH: scroll back to the top (was topLeftCell A10) and move the
# selection to F3 (was A19:D23)
$wsCH = $wb.Worksheets.Item("quotas_CH")
$wsCH.Activate()
$win.ScrollRow = 1
$win.ScrollColumn = 1
$wsCH.Range("F3").Select() | Out-Null

# quotas_MX: move the selection to I18 (was B8:F8)
$wsMX = $wb.Worksheets.Item("quotas_MX")
$wsMX.Activate()
$wsMX.Range("I18").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. Finally activate the last new sheet (quota_UA) with its selection
#    on H4, matching the author's ending view.
# ---------------------------------------------------------------------
$wsUA.Activate()
$wsUA.Range("H4").Select() | Out-Null
